# Sprint_Backlog.xlsx — "alguns acertos na sprint_Backlog"
# Fill in the previously-empty Status (column E) cells on the sprint
# backlog sheet with "Pronto" or "Realizar", and move the active
# selection to E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "Realizar"
$ws.Range("E11").Value = "Realizar"
$ws.Range("E12").Value = "Pronto"
$ws.Range("E15").Value = "Realizar"
$ws.Range("E16").Value = "Realizar"
$ws.Range("E20").Value = "Pronto"
$ws.Range("E26").Value = "Realizar"
$ws.Range("E27").Value = "Pronto"
$ws.Range("E28").Value = "Realizar"
$ws.Range("E30").Value = "Realizar"
$ws.Range("E31").Value = "Realizar"
$ws.Range("E32").Value = "Realizar"
$ws.Range("E33").Value = "Pronto"
$ws.Range("E34").Value = "Pronto"
$ws.Range("E35").Value = "Pronto"

$ws.Range("E13").Select()
